$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("D5").Value = "2016-02-24 07:33:55"
$ws2 = $wb.Worksheets.Item("de-de")
$ws2.Range("D5").Value = "2016-02-24 07:34:07"
